$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.318.04"
$ws.Range("E2").Value = "  +2.83%  "

$ws.Range("D3").Value = "2.018.12"
$ws.Range("E3").Value = "  +6.55%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "246.36"
$ws.Range("E5").Value = "  -0.06%  "

$ws.Range("D6").Value = "0.661"
$ws.Range("E6").Value = "  -4.33%  "

$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "45.14"
$ws.Range("E8").Value = "  +5.41%  "

$ws.Range("D9").Value = "58.82"
$ws.Range("E9").Value = "  +4.44%  "

$ws.Range("D10").Value = "0.361"
$ws.Range("E10").Value = "  +1.03%  "

$ws.Range("D11").Value = "0.0715"
$ws.Range("E11").Value = "  -4.54%  "

$ws.Range("D12").Value = "0.0985"
$ws.Range("E12").Value = "  +0.38%  "

$ws.Range("D13").Value = "14.71"
$ws.Range("E13").Value = "  +5.49%  "

$ws.Range("D14").Value = "2.306.55"
$ws.Range("E14").Value = "  +6.23%  "

$ws.Range("D15").Value = "0.804"
$ws.Range("E15").Value = "  +1.10%  "

$ws.Range("D16").Value = "2.009.66"
$ws.Range("E16").Value = "  +5.64%  "

$ws.Range("D17").Value = "4.91"
$ws.Range("E17").Value = "  -1.98%  "

$ws.Range("D18").Value = "36.025.78"
$ws.Range("E18").Value = "  +1.72%  "

$ws.Range("D19").Value = "71.04"
$ws.Range("E19").Value = "  -3.32%  "

$ws.Range("D20").Value = "0.0₃0819"
$ws.Range("E20").Value = "  -1.08%  "

$ws.Range("D21").Value = "12.95"
$ws.Range("E21").Value = "  +0.15%  "

$ws.Range("D22").Value = "235.30"
$ws.Range("E22").Value = "  -3.66%  "

$ws.Range("D23").Value = "4.89"
$ws.Range("E23").Value = "  -6.30%  "

$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D26").Value = "162.56"
$ws.Range("E26").Value = "  -2.56%  "

$ws.Range("D27").Value = "19.78"
$ws.Range("E27").Value = "  +7.92%  "

$ws.Range("E28").Value = "  -9.43%  "

$ws.Range("D29").Value = "8.56"
$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("D30").Value = "0.122"
$ws.Range("E30").Value = "  -4.36%  "

$ws.Range("D31").Value = "4.38"
$ws.Range("E31").Value = "  +0.65%  "

$ws.Range("D32").Value = "21.15"
$ws.Range("E32").Value = "  +56.33%  "

$ws.Range("D33").Value = "0.0590"
$ws.Range("E33").Value = "  -2.91%  "

$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("E35").Value = "  +0.55%  "

$ws.Range("D36").Value = "4.01"
$ws.Range("E36").Value = "  -5.17%  "

$ws.Range("D37").Value = "0.0801"
$ws.Range("E37").Value = "  +11.36%  "

$ws.Range("D38").Value = "2.12"
$ws.Range("E38").Value = "  +8.07%  "

$ws.Range("D39").Value = "0.841"
$ws.Range("E39").Value = "  -1.28%  "

$ws.Range("D40").Value = "1.35"
$ws.Range("E40").Value = "  -8.31%  "

$ws.Range("D41").Value = "0.0215"
$ws.Range("E41").Value = "  -4.35%  "

$ws.Range("D42").Value = "96.24"
$ws.Range("E42").Value = "  -2.68%  "

$ws.Range("D43").Value = "1.09"
$ws.Range("E43").Value = "  +0.61%  "

$ws.Range("E44").Value = "  +13.76%  "

$ws.Range("D45").Value = "16.02"
$ws.Range("E45").Value = "  -5.65%  "

$ws.Range("D46").Value = "1.317.14"
$ws.Range("E46").Value = "  -1.58%  "

$ws.Range("D47").Value = "0.0809"
$ws.Range("E47").Value = "  -0.04%  "

$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D48").Value = "2.77"
$ws.Range("E48").Value = "  +1.27%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.204.14"
$ws.Range("E49").Value = "  +6.36%  "

$ws.Range("E50").Value = "  -7.26%  "

$ws.Range("D51").Value = "3.79"
$ws.Range("E51").Value = "  +13.39%  "
